$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2271186440677966
$ws.Range("C2").Value = 0.5423728813559322
$ws.Range("J2").Value = 0.01016949152542373
$ws.Range("P2").Value = 0.1491525423728814
$ws.Range("S2").Value = 0.0711864406779661
$ws.Range("B3").Value = 0.005882352941176471
$ws.Range("C3").Value = 0.03529411764705882
$ws.Range("J3").Value = 0.01176470588235294
$ws.Range("P3").Value = 0.7529411764705882
$ws.Range("S3").Value = 0.1941176470588235
$ws.Range("J4").Value = 0.08108108108108109
$ws.Range("P4").Value = 0.7297297297297297
$ws.Range("S4").Value = 0.1891891891891892
$ws.Range("B6").Value = 0.08455882352941177
$ws.Range("F6").Value = 0.06985294117647059
$ws.Range("J6").Value = 0.3014705882352941
$ws.Range("O6").Value = 0.02573529411764706
$ws.Range("Q6").Value = 0.125
$ws.Range("R6").Value = 0.04779411764705882
$ws.Range("S6").Value = 0.3455882352941176
$ws.Range("B7").Value = 0.1057692307692308
$ws.Range("D7").Value = 0.01923076923076923
$ws.Range("F7").Value = 0.0673076923076923
$ws.Range("J7").Value = 0.1201923076923077
$ws.Range("O7").Value = 0.01442307692307692
$ws.Range("Q7").Value = 0.1682692307692308
$ws.Range("R7").Value = 0.07211538461538461
$ws.Range("S7").Value = 0.4326923076923077
$ws.Range("B8").Value = 0.08869179600886919
$ws.Range("D8").Value = 0.01330376940133038
$ws.Range("E8").Value = 0.002217294900221729
$ws.Range("F8").Value = 0.05321507760532151
$ws.Range("J8").Value = 0.1108647450110865
$ws.Range("O8").Value = 0.01995565410199556
$ws.Range("Q8").Value = 0.1596452328159645
$ws.Range("R8").Value = 0.130820399113082
$ws.Range("S8").Value = 0.4212860310421286
$ws.Range("B9").Value = 0.08658008658008658
$ws.Range("D9").Value = 0.02164502164502164
$ws.Range("F9").Value = 0.07792207792207792
$ws.Range("J9").Value = 0.1471861471861472
$ws.Range("O9").Value = 0.02597402597402598
$ws.Range("Q9").Value = 0.1645021645021645
$ws.Range("R9").Value = 0.07792207792207792
$ws.Range("S9").Value = 0.3982683982683983
$ws.Range("B10").Value = 0.1002527379949452
$ws.Range("D10").Value = 0.01853411962931761
$ws.Range("E10").Value = 0.0008424599831508003
$ws.Range("F10").Value = 0.08508845829823083
$ws.Range("J10").Value = 0.1019376579612468
$ws.Range("O10").Value = 0.02274641954507161
$ws.Range("Q10").Value = 0.2080876158382477
$ws.Range("R10").Value = 0.09014321819713564
$ws.Range("S10").Value = 0.3723673125526538
$ws.Range("G11").Value = 0.1155115511551155
$ws.Range("J11").Value = 0.07590759075907591
$ws.Range("K11").Value = 0.1551155115511551
$ws.Range("L11").Value = 0.6336633663366337
$ws.Range("S11").Value = 0.0198019801980198
$ws.Range("G12").Value = 0.765
$ws.Range("J12").Value = 0.18
$ws.Range("K12").Value = 0.005
$ws.Range("L12").Value = 0.03
$ws.Range("S12").Value = 0.02
$ws.Range("G13").Value = 0.6511627906976745
$ws.Range("J13").Value = 0.2558139534883721
$ws.Range("S13").Value = 0.09302325581395349
$ws.Range("F15").Value = 0.04230769230769231
$ws.Range("H15").Value = 0.1615384615384615
$ws.Range("I15").Value = 0.05
$ws.Range("J15").Value = 0.3269230769230769
$ws.Range("K15").Value = 0.08846153846153847
$ws.Range("M15").Value = 0.007692307692307693
$ws.Range("O15").Value = 0.09615384615384616
$ws.Range("S15").Value = 0.2269230769230769
$ws.Range("F16").Value = 0.03157894736842105
$ws.Range("H16").Value = 0.1736842105263158
$ws.Range("I16").Value = 0.1
$ws.Range("J16").Value = 0.3526315789473684
$ws.Range("K16").Value = 0.1052631578947368
$ws.Range("M16").Value = 0.005263157894736842
$ws.Range("O16").Value = 0.07894736842105263
$ws.Range("S16").Value = 0.1526315789473684
$ws.Range("F17").Value = 0.03818615751789976
$ws.Range("H17").Value = 0.1885441527446301
$ws.Range("I17").Value = 0.1097852028639618
$ws.Range("J17").Value = 0.3890214797136038
$ws.Range("K17").Value = 0.1026252983293556
$ws.Range("M17").Value = 0.02863961813842482
$ws.Range("O17").Value = 0.07875894988066826
$ws.Range("S17").Value = 0.06443914081145585
$ws.Range("F18").Value = 0.02790697674418605
$ws.Range("H18").Value = 0.1674418604651163
$ws.Range("I18").Value = 0.1209302325581395
$ws.Range("J18").Value = 0.386046511627907
$ws.Range("K18").Value = 0.1116279069767442
$ws.Range("M18").Value = 0.01395348837209302
$ws.Range("O18").Value = 0.09302325581395349
$ws.Range("S18").Value = 0.07906976744186046
$ws.Range("F19").Value = 0.02780049059689289
$ws.Range("H19").Value = 0.2150449713818479
$ws.Range("I19").Value = 0.1046606704824203
$ws.Range("J19").Value = 0.3401471790678659
$ws.Range("K19").Value = 0.116107931316435
$ws.Range("M19").Value = 0.02289452166802943
$ws.Range("N19").Value = 0.001635322976287817
$ws.Range("O19").Value = 0.07031888798037612
$ws.Range("S19").Value = 0.1013900245298446
